# Apply the commit "Add files via upload" changes to the workbook.
# This fills in the previously-empty row 13 (19/04/2017) on sheet "Abril"
# with a time entry and description, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Abril")

# Entrada / Saída times for row 13 (8:00 and 10:00), stored as the
# fraction-of-a-day values Excel uses internally for time-only cells.
$ws.Range("B13").Value = 0.33333333333333331
$ws.Range("C13").Value = 0.41666666666666669

# Atividade Desenvolvida / Descrição text for row 13
$ws.Range("E13").Value = "Modelagem com Blender"
$ws.Range("F13").Value = "Pesquisa sobre outras formas de rotacionar as faces dos poliedros para planificá-los"

# Update the selection to match the edited cell, like a user just finished
# typing into F13.
$ws.Range("F13").Select()

$wb.Save()
